# Improve the Lexical Analyzer sheet: remove the "Code" column (column G)
# from the instruction tables on Sheet1. This shifts the "Operandos" column
# (and the hex value columns in the data rows) one column to the left, and
# removes the now-unused "Code"/"0x60" shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire "Code" column (G). Excel will automatically shift all
# columns to its right (Operandos, etc.) one position to the left and will
# drop shared-string entries that are no longer referenced by any cell.
$ws.Columns("G").Delete() | Out-Null

# Reflect where the user ended up clicking after performing the edit.
$ws.Range("E7").Select() | Out-Null
